# Update the "EPS based on latest capital" row (row 27, columns D:H)
# to reflect the revised read_price algorithm / refreshed database values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("D27").Value = 77
$ws.Range("E27").Value = 335
$ws.Range("F27").Value = 306
$ws.Range("G27").Value = 1062
$ws.Range("H27").Value = 1637
